# Refresh the Khanty-Mansiysk forecast sheet with the latest model run.
# Columns: A=Number(units) B=Area(ha) C=Forest area(ha)
#          H=Forecast Forest area(ha) I=Forecast Area(ha) J=Forecast Number(units)
# D/E/F/G (Year, Accumulated temperature/precipitations) are untouched inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @(2,  300,  76657,  9210,  -56647,  -57179,  99),
    @(3,  46,   2868,   646,   -89845,  -105741, 320),
    @(4,  126,  6710,   898,   31251,   39994,   135),
    @(5,  720,  119922, 9564,  19178,   56900,   571),
    @(6,  959,  276787, 34271, 136886,  282543,  708),
    @(7,  912,  198551, 58652, 196804,  377611,  759),
    @(8,  453,  59529,  10614, 6904,    48306,   523),
    @(9,  336,  47835,  26795, -25653,  -11264,  444),
    @(10, 237,  23477,  5100,  -17228,  -8563,   424),
    @(11, 529,  69644,  8732,  10476,   39825,   495),
    @(12, 415,  107979, 61888, 47479,   80907,   517),
    @(13, 825,  188100, 93465, 148930,  286757,  701),
    @(14, 1235, 1317798,786198,686201,  1163573, 1045),
    @(15, 517,  200479, 134647,29505,   87415,   561),
    @(16, 94,   6074,   1775,  66922,   96754,   64),
    @(17, 122,  10418,  2338,  39391,   96187,   -46),
    @(18, 244,  28917,  18388, 39299,   66285,   577),
    @(19, 383,  99483,  81724, 64934,   154783,  621),
    @(20, 558,  24722,  12674, 1597,    -3711,   361),
    @(21, 543,  31199,  9198,  17716,   76254,   562),
    @(22, 633,  227985, 185155,197834,  357497,  747)
)

foreach ($row in $updates) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value  = $row[1]   # A - Number (units)
    $ws.Cells.Item($r, 2).Value  = $row[2]   # B - Area (ha)
    $ws.Cells.Item($r, 3).Value  = $row[3]   # C - Forest area (ha)
    $ws.Cells.Item($r, 8).Value  = $row[4]   # H - Forecast Forest area (ha)
    $ws.Cells.Item($r, 9).Value  = $row[5]   # I - Forecast Area (ha)
    $ws.Cells.Item($r, 10).Value = $row[6]   # J - Forecast Number (units)
}

# Column C got a bit narrower now that it only needs to fit integers.
$ws.Columns.Item(3).ColumnWidth = 20.83
